$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark rows 23-27 (Labyrinthine through Paradigm) as "Processed" in column C,
# matching the existing pattern used for rows 2-22.
for ($r = 23; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "Processed"
}
